$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.904.33"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").Value = "3.224.52"
$ws.Range("E3").Value = "  -4.11%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'539.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.74%  "
$ws.Range("D6").Value = "'136.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.76%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.223.34"
$ws.Range("E8").Value = "  -4.18%  "
$ws.Range("E9").Value = "  -4.35%  "
$ws.Range("D10").Value = "'7.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.05%  "
$ws.Range("E11").Value = "  -6.02%  "
$ws.Range("D12").Value = "'0.395"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.54%  "
$ws.Range("D13").Value = "3.779.71"
$ws.Range("E13").Value = "  -4.01%  "
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'26.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.32%  "
$ws.Range("D16").Value = "3.231.37"
$ws.Range("E16").Value = "  -3.67%  "
$ws.Range("E17").Value = "  -6.57%  "
$ws.Range("D18").Value = "58.952.93"
$ws.Range("E18").Value = "  -3.45%  "
$ws.Range("E19").Value = "  -6.46%  "
$ws.Range("D20").Value = "'13.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.31%  "
$ws.Range("E21").Value = "  -6.87%  "
$ws.Range("D22").Value = "'363.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'70.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.14%  "
$ws.Range("E25").Value = "  -6.97%  "
$ws.Range("D26").Value = "3.360.45"
$ws.Range("E26").Value = "  -4.00%  "
$ws.Range("D27").Value = "'0.172"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").Value = "0.0₃0967"
$ws.Range("E28").Value = "  -10.91%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'7.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.20%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "'1.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.77%  "
$ws.Range("E33").Value = "  -8.23%  "
$ws.Range("D34").Value = "'21.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("E35").Value = "  -7.01%  "
$ws.Range("E36").Value = "  -8.38%  "
$ws.Range("D37").Value = "'161.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.33%  "
$ws.Range("E38").Value = "  -5.59%  "
$ws.Range("E39").Value = "  -6.56%  "
$ws.Range("D40").Value = "'26.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.69%  "
$ws.Range("D41").Value = "'0.0709"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.86%  "
$ws.Range("D42").Value = "3.254.01"
$ws.Range("E42").Value = "  -4.19%  "
$ws.Range("D43").Value = "'41.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("E44").Value = "  -5.68%  "
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("E46").Value = "  -6.14%  "
$ws.Range("D47").Value = "'1.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.60%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "2.304.24"
$ws.Range("E49").Value = "  -7.38%  "
$ws.Range("D50").Value = "'6.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.37%  "
$ws.Range("D51").Value = "'20.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.33%  "
